# Apply update described by the diff:
# - Column C (Förändrad) rows 2..34: 45644 -> 45645
# - Row 34 loses its custom row height (ht=15 customHeight=1)
# - Rows 35 and 36 are deleted
# - dimension shrinks from A1:Z36 to A1:Z34

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values for rows 2 through 34 (the "Förändrad" date column).
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = 45645
}

# Delete rows 35 and 36 entirely.
$ws.Rows.Item(36).Delete()
$ws.Rows.Item(35).Delete()

# Row 34 should no longer have an explicit custom row height (ht/customHeight
# attributes removed) - AutoFit drops the custom height override.
$ws.Rows.Item(34).AutoFit()
